$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 352
$ws1.Range("F7").Value = 1156
$ws1.Range("F9").Value = 7078
$ws1.Range("F13").Value = 7949
$ws1.Range("F16").Value = 5497
$ws1.Range("F18").Value = 2391
$ws1.Range("F20").Value = 4558
$ws1.Range("F23").Value = 77
$ws1.Range("F25").Value = 363
$ws1.Range("F28").Value = 2279
$ws1.Range("F30").Value = 260
$ws1.Range("F31").Value = 71
$ws1.Range("F32").Value = 125
$ws1.Range("F33").Value = 572
$ws1.Range("F36").Value = 1473
$ws1.Range("F39").Value = 2275
$ws1.Range("F40").Value = 2207

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 58
$ws2.Range("F8").Value = 94

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1274
$ws4.Range("F7").Value = 352
$ws4.Range("F8").Value = 1156
$ws4.Range("F10").Value = 7078
$ws4.Range("F14").Value = 7949
$ws4.Range("F17").Value = 5497
$ws4.Range("F19").Value = 2391
$ws4.Range("F21").Value = 4558
$ws4.Range("F24").Value = 77
$ws4.Range("F27").Value = 58
$ws4.Range("F28").Value = 363
$ws4.Range("F30").Value = 2279
$ws4.Range("F32").Value = 260
$ws4.Range("F33").Value = 71
$ws4.Range("F34").Value = 125
$ws4.Range("F36").Value = 572
$ws4.Range("F40").Value = 1473
$ws4.Range("F43").Value = 2275
$ws4.Range("F45").Value = 2207
$ws4.Range("F49").Value = 94
